$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraper re-ran and re-appended its last 10 rows (944:953) to the end
# of the sheet (954:963) verbatim - duplicate rows with updated row numbers.
$ws.Range("A944:W953").Copy() | Out-Null
$excel.ActiveSheet.Paste($ws.Range("A954")) | Out-Null

# Re-apply the "best odds" red highlight (fill) to the copied max-value
# cells - copy/paste above only carried over values, not formatting.
$ws.Range("I954").Interior.Color = 255
$ws.Range("I956").Interior.Color = 255
$ws.Range("I957").Interior.Color = 255
$ws.Range("I958").Interior.Color = 255
$ws.Range("E959").Interior.Color = 255
$ws.Range("I959").Interior.Color = 255
$ws.Range("E961").Interior.Color = 255
$ws.Range("Q961").Interior.Color = 255
$ws.Range("U961").Interior.Color = 255
$ws.Range("I962").Interior.Color = 255
$ws.Range("I963").Interior.Color = 255

# Move the selection to the new last row, one past the appended data,
# matching the author's cursor position after the paste.
$ws.Range("A964").Select() | Out-Null
